# Revision del mes de julio
# Hides previously-closed rows (now filtered out by the "Abierto" status
# AutoFilter), records the closing date + new status for row 15, flags
# row 16 as a risk ("Si"), and logs a new follow-up action in row 17
# (copy of row 16's deviation, opened in August).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: real closing date + status -> Cerrado ---
$ws.Cells.Item(15, 10).Value2 = 42593      # J15 Fecha Real Cierre = 11/08/2016
$ws.Cells.Item(15, 11).Value2 = "Cerrado"  # K15 Estatus

# --- Row 16: mark Riesgo = Si ---
$ws.Cells.Item(16, 4).Value2 = "Si"        # D16 Riesgo

# --- Row 17: new follow-up entry (same deviation as row 16, new dates) ---
$ws.Cells.Item(17, 4).Value2 = "Si"
$ws.Cells.Item(17, 5).Value2 = "Presentar esfuerzos superiores a los estimados en la cotización anual de la empresa lo cual provoca perdidas monetarias a la empresa"
$ws.Cells.Item(17, 6).Value2 = "En espera de acciones correctivas"
$ws.Cells.Item(17, 7).Value2 = 42593       # G17 Fecha Compromiso = 11/08/2016
$ws.Cells.Item(17, 8).Value2 = "Ricardo Novela"
$ws.Cells.Item(17, 9).Value2 = 42597       # I17 Fecha de Deteccion = 15/08/2016
$ws.Cells.Item(17, 11).Value2 = "Abierto"

# Match J17's number-format to the other date cells in the row (copy format
# only, leave the cell itself blank, same as the source diff).
$ws.Cells.Item(16, 10).Copy()
$ws.Cells.Item(17, 10).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 17 now holds wrapped, multi-line text like rows 14/16, so it grows to
# the same row height.
$ws.Rows.Item(17).RowHeight = 114

# --- Hide rows 13, 14 and 15 (their "Estatus" is/becomes "Cerrado", which
#     the K-column AutoFilter (Abierto + blanks) no longer shows). Done
#     after the cell writes above so the stored row heights aren't
#     recomputed away by the hidden-row (0-height) state. ---
$ws.Rows.Item(13).Hidden = $true
$ws.Rows.Item(14).Hidden = $true
$ws.Rows.Item(15).Hidden = $true

# Cursor ends up parked on the freshly edited row.
$ws.Range("R17").Select()
